# Append a new data row (row 24) continuing the existing "CodeConfig" table,
# mirroring the pattern of the preceding rows (C=#, D=code, E=type, F=id, G=qty).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C24").Value = 19
$ws.Range("D24").Value = "gx142"
$ws.Range("E24").Value = 9
$ws.Range("F24").Value = 4003
$ws.Range("G24").Value = 100

# Mirror the saved view/selection state (activeCell moved from G28 to G29).
$ws.Range("G29").Select()
